$wb = $excel.ActiveWorkbook

# --- Core content change: delete row 37 on "BTS NTS Modal Profile Data" ---
# Row 37 held the "weighted value, adjusted for number of train cars per
# locomotive" helper line (A37/B37 = B36/10). Removing the row shifts every
# row below it up by one.
$modal = $wb.Worksheets.Item("BTS NTS Modal Profile Data")
$modal.Rows.Item(37).Delete()

# --- Fix up cross-sheet references that pointed at the now-deleted row ---
# Row-shift within the same workbook normally re-targets formulas
# automatically, but the formula on "AVLo-passengers" used to reference the
# now-removed B37 (the weighted/derived value) and must instead point at the
# original B36 value (the un-derived Average Vehicle Loading figure).
$passengers = $wb.Worksheets.Item("AVLo-passengers")
$passengers.Range("B5").Formula = "='BTS NTS Modal Profile Data'!B36"

# --- View-state cleanup to mirror the reverted workbook ---
# The "reverted" file was last saved with the "About" sheet active/selected
# (cell A44) instead of "AVLo-passengers" (which had been the active tab,
# scrolled/selected at F12). Re-home the selections accordingly.
$null = $modal.Select()
$null = $modal.Range("A1").Select()

$null = $passengers.Select()
$null = $passengers.Range("A1").Select()

$about = $wb.Worksheets.Item("About")
$null = $about.Select()
$null = $about.Range("A44").Select()

$wb.Application.CalculateFull() | Out-Null
